# Apply the "Saldo" worksheet updates.
#
# The sheet has columns: A=Conta (account, stored as text w/ leading zeros),
# B=Nome (name, text), C=Saldo (numeric balance). Row 1 is the header row,
# so data row N in the sheet is spreadsheet row N (1-indexed, header = row 1).
#
# Operations are applied from the bottom of the sheet upward so that row
# numbers referenced below are never invalidated by an earlier insert/delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# 1) Row 336: account 004212581 (MARIA, -29543.97) is removed from this
#    spot entirely (the account re-appears higher up with a new balance,
#    handled in step 7 below).
# ---------------------------------------------------------------------
$ws.Rows(336).Delete()

# ---------------------------------------------------------------------
# 2) Row 221: account 004207184 (CRISTINA) balance changes
#    from -1960.26 to -1911.33.
# ---------------------------------------------------------------------
$ws.Range("C221").Value = -1911.33

# ---------------------------------------------------------------------
# 3) Row 49: old account 004205505 (SURAMA, 512.08) row is removed from
#    this spot (the account re-appears a bit higher with a new balance,
#    handled in step 4 below).
# ---------------------------------------------------------------------
$ws.Rows(49).Delete()

# ---------------------------------------------------------------------
# 4) Insert a new row above row 32 (005044389 CLAUDIA) for account
#    004205505 (SURAMA) with its new balance of 756.72.
# ---------------------------------------------------------------------
$ws.Rows(32).Insert()
Set-TextCell $ws.Range("A32") "004205505"
Set-TextCell $ws.Range("B32") "SURAMA"
$ws.Range("C32").Value = 756.72

# ---------------------------------------------------------------------
# 5) Row 9: account 004212476 (MARIA) balance changes
#    from 50247.26 to 63076.55.
# ---------------------------------------------------------------------
$ws.Range("C9").Value = 63076.55

# ---------------------------------------------------------------------
# 6) Insert a new row right after row 9 (i.e. above the current row 10)
#    for a brand-new account 004328934 (VALERIA) with balance 30000.
# ---------------------------------------------------------------------
$ws.Rows(10).Insert()
Set-TextCell $ws.Range("A10") "004328934"
Set-TextCell $ws.Range("B10") "VALERIA"
$ws.Range("C10").Value = 30000

# ---------------------------------------------------------------------
# 7) Insert a new row above row 7 (004218542 JOSE) for account
#    004212581 (MARIA) with its new balance of 93601.44.
# ---------------------------------------------------------------------
$ws.Rows(7).Insert()
Set-TextCell $ws.Range("A7") "004212581"
Set-TextCell $ws.Range("B7") "MARIA"
$ws.Range("C7").Value = 93601.44
